$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 12
$ws.Cells.Item(2, 2).Value = 100
$ws.Cells.Item(2, 3).Value = 14
$ws.Cells.Item(2, 4).Value = 198
$ws.Cells.Item(2, 5).Value = 184
$ws.Cells.Item(2, 6).Value = 3
$ws.Cells.Item(2, 7).Value = 1
$ws.Cells.Item(2, 8).Value = 1
$ws.Cells.Item(2, 9).Value = 1100

$ws.Cells.Item(3, 1).Value = 4
$ws.Cells.Item(3, 2).Value = 7
$ws.Cells.Item(3, 3).Value = 14
$ws.Cells.Item(3, 4).Value = 195
$ws.Cells.Item(3, 5).Value = 179
$ws.Cells.Item(3, 6).Value = 3
$ws.Cells.Item(3, 7).Value = 1
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 91
